$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The helper lookup table that lived in column M (header "Communes" in M1,
# looked-up commune names in M3:M11) is no longer needed, so the whole
# column is removed.
$ws.Columns("M").Delete()

# Rows 39-43 (the "Ratissage" period) didn't have a "Lots" label yet; fill
# it in like the surrounding rows already do for their own periods.
$ws.Range("E39:E43").Value = "Ratissage"

# Restore the view state (scrolled down a bit, selection moved along).
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("G41").Select() | Out-Null
